$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row: Fecha(D), Volumen(J), Precio minimo(K), Precio maximo(L), Precio promedio ponderado(M), Precio $/Kg(P)
$rows = @{
    2  = @(44320, 160, 19000, 20000, 19500, 1500)
    3  = @(44580, 160, 11000, 12000, 11500, 885)
    4  = @(44764, 200, 12000, 13000, 12500, 962)
    5  = @(44616, 120, 19000, 20000, 19500, 1500)
    6  = @(44397, 140, 12500, 13000, 12750, 981)
    7  = @(44159, 100, 23000, 24000, 23500, 1808)
    8  = @(44469, 140, 13000, 14000, 13500, 1038)
    9  = @(44389, 120, 12000, 13000, 12500, 962)
    10 = @(44592, 120, 12000, 13000, 12500, 962)
    11 = @(44379, 120, 12000, 13000, 12667, 974)
    12 = @(44832, 100, 13000, 14000, 13500, 1038)
    13 = @(44406, 160, 17000, 18000, 17500, 1346)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals[0]
    $ws.Range("J$r").Value = $vals[1]
    $ws.Range("K$r").Value = $vals[2]
    $ws.Range("L$r").Value = $vals[3]
    $ws.Range("M$r").Value = $vals[4]
    $ws.Range("P$r").Value = $vals[5]
}
